# Update Name of Algo
# Applies the updated imputed values produced by the KNN algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B4"  = 6.821000000000001
    "A9"  = -21.846
    "B9"  = 5.518000000000001
    "C9"  = -11.169
    "B11" = 6.108
    "A13" = -21.844
    "A16" = -20.879
    "B16" = 6.677000000000001
    "A18" = -21.798
    "A20" = -20.027
    "C22" = -12.329
    "B23" = 6.842000000000001
    "B24" = 6.220000000000001
    "A26" = -20.959
    "B26" = 7.002
    "A27" = -21.495
    "C27" = -13.017
    "A29" = -21.333
    "C29" = -12.053
    "C32" = -12.633
    "B34" = 7.256
    "A35" = -21.37
    "B35" = 6.333
    "A36" = -20.714
    "C37" = -12.086
    "C38" = -11.714
    "C39" = -12.741
    "C41" = -12.546
    "B44" = 5.790000000000001
    "A45" = -21.422
    "C45" = -12.819
    "B48" = 5.518000000000001
    "C48" = -11.328
    "B49" = 5.93
    "C51" = -11.235
    "B52" = 5.269
    "A55" = -22.124
    "C56" = -12.5
    "A57" = -21.646
    "C57" = -13.216
    "C61" = -12.624
    "C64" = -11.023
    "B66" = 5.491
    "B67" = 5.012
    "A69" = -21.361
    "B73" = 5.745
    "C75" = -12.43
    "A76" = -20.293
    "A78" = -20.923
    "B78" = 6.795
    "B80" = 6.738000000000001
    "A82" = -21.258
    "C82" = -11.541
    "A83" = -21.529
    "C90" = -11.009
    "B91" = 5.469
    "A93" = -21.305
    "C93" = -11.708
    "A97" = -21.46100000000001
    "B97" = 5.3
    "B99" = 5.347
    "C102" = -12.456
    "B104" = 7.495
    "C105" = -12.498
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
